$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold an exact text value (prevents Excel from
# auto-converting number-looking strings such as "23.64" or "1.009" into
# floating point numbers, which would corrupt the stored text).
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
}

# Row 2
Set-TextValue $ws.Range("D2") '29.554.88'
$ws.Range("E2").Value = '  +0.09%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.921.79'
$ws.Range("E3").Value = '  +0.35%  '

# Row 4
$ws.Range("E4").Value = '  +0.39%  '

# Row 5
Set-TextValue $ws.Range("D5") '325.78'
$ws.Range("E5").Value = '  -0.01%  '

# Row 6
Set-TextValue $ws.Range("D6") '1.009'
$ws.Range("E6").Value = '  +0.37%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.4814'
$ws.Range("E7").Value = '  -0.17%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.4069'
$ws.Range("E8").Value = '  -0.16%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.08224'
$ws.Range("E9").Value = '  +0.85%  '

# Row 10
Set-TextValue $ws.Range("D10") '1.010'
$ws.Range("E10").Value = '  -0.32%  '

# Row 11
Set-TextValue $ws.Range("D11") '23.64'
$ws.Range("E11").Value = '  +1.05%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D12") '6.081'
$ws.Range("E12").Value = '  +1.25%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D13") '1.895.39'
$ws.Range("E13").Value = '  -1.25%  '

# Row 14
Set-TextValue $ws.Range("D14") '7.267'
$ws.Range("E14").Value = '  +1.96%  '

# Row 15
Set-TextValue $ws.Range("D15") '91.67'
$ws.Range("E15").Value = '  +1.52%  '

# Row 16
Set-TextValue $ws.Range("D16") '0.06863'

# Row 17
$ws.Range("E17").Value = '  +0.32%  '

# Row 18
$ws.Range("E18").Value = '  -0.11%  '

# Row 19
Set-TextValue $ws.Range("D19") '17.64'
$ws.Range("E19").Value = '  -0.30%  '

# Row 20
Set-TextValue $ws.Range("D20") '1.010'
$ws.Range("E20").Value = '  +0.29%  '

# Row 21
Set-TextValue $ws.Range("D21") '29.564.51'
$ws.Range("E21").Value = '  +0.05%  '

# Row 22
Set-TextValue $ws.Range("D22") '5.681'
$ws.Range("E22").Value = '  +1.07%  '

# Row 23
Set-TextValue $ws.Range("D23") '11.93'
$ws.Range("E23").Value = '  +1.11%  '

# Row 24
Set-TextValue $ws.Range("D24") '2.183'
$ws.Range("E24").Value = '  +0.05%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.128.60'
$ws.Range("E25").Value = '  -0.88%  '

# Row 26
Set-TextValue $ws.Range("D26") '155.96'
$ws.Range("E26").Value = '  +0.42%  '

# Row 27
Set-TextValue $ws.Range("D27") '6.444'
$ws.Range("E27").Value = '  -1.03%  '

# Row 28
Set-TextValue $ws.Range("D28") '20.01'
$ws.Range("E28").Value = '  -0.12%  '

# Row 29
Set-TextValue $ws.Range("D29") '2.096'
$ws.Range("E29").Value = '  -0.16%  '

# Row 30
Set-TextValue $ws.Range("D30") '120.60'
$ws.Range("E30").Value = '  +0.77%  '

# Row 31
Set-TextValue $ws.Range("D31") '1.014'
$ws.Range("E31").Value = '  -1.70%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.09612'
$ws.Range("E32").Value = '  +0.44%  '

# Row 33
Set-TextValue $ws.Range("D33") '5.621'
$ws.Range("E33").Value = '  +1.87%  '

# Row 34
Set-TextValue $ws.Range("D34") '3.550'
$ws.Range("E34").Value = '  -0.29%  '

# Row 35
$ws.Range("E35").Value = '  -1.18%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.06358'
$ws.Range("E36").Value = '  +4.17%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.02293'
$ws.Range("E37").Value = '  +1.08%  '

# Row 38
Set-TextValue $ws.Range("D38") '1.185'
$ws.Range("E38").Value = '  +0.44%  '

# Row 39
Set-TextValue $ws.Range("D39") '0.5958'
$ws.Range("E39").Value = '  +0.39%  '

# Row 40
Set-TextValue $ws.Range("D40") '10.78'
$ws.Range("E40").Value = '  +0.05%  '

# Row 41
$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D41") '1.009'
$ws.Range("E41").Value = '  +0.30%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D42") '7.884'
$ws.Range("E42").Value = '  -1.06%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.1850'
$ws.Range("E43").Value = '  -0.45%  '

# Row 44
Set-TextValue $ws.Range("D44") '2.456'
$ws.Range("E44").Value = '  -1.22%  '

# Row 45
Set-TextValue $ws.Range("D45") '1.270'
$ws.Range("E45").Value = '  -1.03%  '

# Row 46
Set-TextValue $ws.Range("D46") '12.45'
$ws.Range("E46").Value = '  -0.39%  '

# Row 47
Set-TextValue $ws.Range("D47") '0.07490'
$ws.Range("E47").Value = '  -2.96%  '

# Row 48
Set-TextValue $ws.Range("D48") '0.5563'
$ws.Range("E48").Value = '  -0.11%  '

# Row 50
Set-TextValue $ws.Range("D50") '119.30'
$ws.Range("E50").Value = '  +3.22%  '

# Row 51
$ws.Range("E51").Value = '  +3.17%  '
